# The deck's cached "today" date placeholders (Date Placeholder shapes on the
# Slide Master, every Slide Layout, the Notes Master and the Handout Master)
# were refreshed from 6/5/2017 to 6/6/2017 (e.g. the presentation was opened
# and re-saved a day later). Walk every one of those tiers and update the
# cached date text wherever it still shows the old date.

function Update-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "6/5/2017") {
                $tr.Text = "6/6/2017"
            }
        }
    }
}

$p = $ppt.ActivePresentation

# Slide Master
$master = $p.SlideMaster
Update-DateText $master.Shapes

# Every Slide Layout hanging off the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateText $layout.Shapes
}

# Notes Master
if ($p.HasNotesMaster) {
    Update-DateText $p.NotesMaster.Shapes
}

# Handout Master
if ($p.HasHandoutMaster) {
    Update-DateText $p.HandoutMaster.Shapes
}
